# Update cryptocurrency price/volume figures ("cryptos list" refresh).
# D column values are apostrophe-prefixed so Excel stores them as text
# (preserving formats like "27.912.17" / "213.52" / "1.80") instead of
# silently coercing them to floating point numbers; the style is reset
# back to "Normal" afterwards so no stray number-format/quote-prefix
# style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.912.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = "'1.642.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'213.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'23.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.50%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").Value = "'0.0874"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = "'1.875.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("D13").Value = "'1.645.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = "'0.574"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.30%  '
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").Value = "'65.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = "'27.891.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = "'230.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = "'10.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.02%  '
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("E24").Value = '  +2.25%  '
$ws.Range("D25").Value = "'152.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("D26").Value = "'6.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").Value = "'15.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("D33").Value = "'1.426.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.60%  '
$ws.Range("D34").Value = "'3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("E35").Value = '  +1.55%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").Value = "'0.889"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.86%  '
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("D39").Value = "'0.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("D40").Value = "'0.557"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").Value = "'1.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = "'68.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("D45").Value = "'5.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.00%  '
$ws.Range("D46").Value = "'1.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").Value = "'1.784.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").Value = "'88.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.83%  '
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("E51").Value = '  +0.53%  '
